$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Cells.Item(5, 2).Value = 6726054
$ws.Cells.Item(5, 7).Value = "Shakhter Karagandy"
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 1.727
$ws.Cells.Item(5, 12).Value = 4
$ws.Cells.Item(5, 13).Value = 3.5
$ws.Cells.Item(5, 14).Value = 1.333
$ws.Cells.Item(5, 15).Value = 4.75
$ws.Cells.Item(5, 16).Value = 6.5
$ws.Cells.Item(5, 17).Value = -1.5
$ws.Cells.Item(5, 18).Value = 1.975
$ws.Cells.Item(5, 19).Value = 1.825
$ws.Cells.Item(5, 20).Value = 3
$ws.Cells.Item(5, 21).Value = 1.975
$ws.Cells.Item(5, 22).Value = 1.825
$ws.Cells.Item(5, 23).Value = 0.333
$ws.Cells.Item(5, 26).Value = 0.9750000000000001
$ws.Cells.Item(5, 27).Value = -1
$ws.Cells.Item(5, 28).Value = -1
$ws.Cells.Item(5, 29).Value = 0.825
# Row 6
$ws.Cells.Item(6, 2).Value = 6221786
$ws.Cells.Item(6, 7).Value = "FK Atyrau"
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 11).Value = 1.571
$ws.Cells.Item(6, 12).Value = 3.8
$ws.Cells.Item(6, 13).Value = 4.75
$ws.Cells.Item(6, 14).Value = 1.5
$ws.Cells.Item(6, 15).Value = 4
$ws.Cells.Item(6, 16).Value = 5.25
$ws.Cells.Item(6, 17).Value = -1
$ws.Cells.Item(6, 18).Value = 1.875
$ws.Cells.Item(6, 19).Value = 1.925
$ws.Cells.Item(6, 20).Value = 2.5
$ws.Cells.Item(6, 21).Value = 1.9
$ws.Cells.Item(6, 22).Value = 1.9
$ws.Cells.Item(6, 23).Value = 0.5
$ws.Cells.Item(6, 26).Value = 0
$ws.Cells.Item(6, 27).Value = -0
$ws.Cells.Item(6, 28).Value = 0.8999999999999999
$ws.Cells.Item(6, 29).Value = -1
# Row 9
$ws.Cells.Item(9, 2).Value = 6221694
$ws.Cells.Item(9, 6).Value = "FC Astana"
$ws.Cells.Item(9, 7).Value = "FK Kaspyi Aktau"
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 11).Value = 1.333
$ws.Cells.Item(9, 12).Value = 4.333
$ws.Cells.Item(9, 13).Value = 7.5
$ws.Cells.Item(9, 14).Value = 1.2
$ws.Cells.Item(9, 15).Value = 5.5
$ws.Cells.Item(9, 16).Value = 11
$ws.Cells.Item(9, 17).Value = -2
$ws.Cells.Item(9, 18).Value = 1.975
$ws.Cells.Item(9, 19).Value = 1.825
$ws.Cells.Item(9, 20).Value = 3
$ws.Cells.Item(9, 21).Value = 1.95
$ws.Cells.Item(9, 22).Value = 1.85
$ws.Cells.Item(9, 23).Value = 0.2
$ws.Cells.Item(9, 26).Value = 0
$ws.Cells.Item(9, 27).Value = -0
$ws.Cells.Item(9, 28).Value = 0.95
# Row 10
$ws.Cells.Item(10, 2).Value = 6221789
$ws.Cells.Item(10, 6).Value = "Kairat Almaty"
$ws.Cells.Item(10, 7).Value = "FK Kyzylzhar"
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 11).Value = 2.25
$ws.Cells.Item(10, 12).Value = 3.2
$ws.Cells.Item(10, 13).Value = 2.8
$ws.Cells.Item(10, 14).Value = 2.1
$ws.Cells.Item(10, 15).Value = 3
$ws.Cells.Item(10, 16).Value = 3.2
$ws.Cells.Item(10, 17).Value = -0.25
$ws.Cells.Item(10, 18).Value = 1.875
$ws.Cells.Item(10, 19).Value = 1.925
$ws.Cells.Item(10, 20).Value = 2.25
$ws.Cells.Item(10, 21).Value = 2
$ws.Cells.Item(10, 22).Value = 1.8
$ws.Cells.Item(10, 23).Value = 1.1
$ws.Cells.Item(10, 26).Value = 0.875
$ws.Cells.Item(10, 27).Value = -1
$ws.Cells.Item(10, 28).Value = 1
# Row 13
$ws.Cells.Item(13, 7).Value = "FK Aktobe"
# Row 15
$ws.Cells.Item(15, 6).Value = "Ordabasy"
# Row 16
$ws.Cells.Item(16, 7).Value = "Ordabasy"
# Row 22
$ws.Cells.Item(22, 6).Value = "FK Aktobe"
# Row 23
$ws.Cells.Item(23, 7).Value = "FK Aktobe"
# Row 25
$ws.Cells.Item(25, 7).Value = "Ordabasy"
# Row 27
$ws.Cells.Item(27, 7).Value = "FK Aktobe"
# Row 31
$ws.Cells.Item(31, 6).Value = "Ordabasy"
# Row 32
$ws.Cells.Item(32, 6).Value = "FK Aktobe"
# Row 39
$ws.Cells.Item(39, 7).Value = "Ordabasy"
# Row 45
$ws.Cells.Item(45, 6).Value = "Ordabasy"
# Row 50
$ws.Cells.Item(50, 2).Value = 6221723
$ws.Cells.Item(50, 6).Value = "FK Aksu"
$ws.Cells.Item(50, 8).Value = 1
$ws.Cells.Item(50, 9).Value = 1
$ws.Cells.Item(50, 10).Value = "D"
$ws.Cells.Item(50, 11).Value = 3.75
$ws.Cells.Item(50, 12).Value = 3.3
$ws.Cells.Item(50, 13).Value = 1.833
$ws.Cells.Item(50, 14).Value = 2.6
$ws.Cells.Item(50, 15).Value = 3.1
$ws.Cells.Item(50, 16).Value = 2.5
$ws.Cells.Item(50, 17).Value = 0
$ws.Cells.Item(50, 18).Value = 1.925
$ws.Cells.Item(50, 19).Value = 1.875
$ws.Cells.Item(50, 20).Value = 2.5
$ws.Cells.Item(50, 21).Value = 1.9
$ws.Cells.Item(50, 22).Value = 1.9
$ws.Cells.Item(50, 24).Value = 2.1
$ws.Cells.Item(50, 25).Value = -1
$ws.Cells.Item(50, 28).Value = -1
$ws.Cells.Item(50, 29).Value = 0.8999999999999999
# Row 51
$ws.Cells.Item(51, 2).Value = 7055064
$ws.Cells.Item(51, 6).Value = "Kaisar Kyzylorda"
$ws.Cells.Item(51, 8).Value = 2
$ws.Cells.Item(51, 9).Value = 3
$ws.Cells.Item(51, 10).Value = "A"
$ws.Cells.Item(51, 11).Value = 6.5
$ws.Cells.Item(51, 12).Value = 4.5
$ws.Cells.Item(51, 13).Value = 1.363
$ws.Cells.Item(51, 14).Value = 4.2
$ws.Cells.Item(51, 15).Value = 4
$ws.Cells.Item(51, 16).Value = 1.6
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = 1.725
$ws.Cells.Item(51, 19).Value = 2.075
$ws.Cells.Item(51, 20).Value = 2.25
$ws.Cells.Item(51, 21).Value = 1.875
$ws.Cells.Item(51, 22).Value = 1.925
$ws.Cells.Item(51, 24).Value = -1
$ws.Cells.Item(51, 25).Value = 0.6000000000000001
$ws.Cells.Item(51, 28).Value = 0.875
$ws.Cells.Item(51, 29).Value = -1
# Row 57
$ws.Cells.Item(57, 6).Value = "Ordabasy"
# Row 58
$ws.Cells.Item(58, 6).Value = "FK Aktobe"
# Row 63
$ws.Cells.Item(63, 2).Value = 6221732
$ws.Cells.Item(63, 6).Value = "FK Atyrau"
$ws.Cells.Item(63, 7).Value = "Ordabasy"
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 10).Value = "D"
$ws.Cells.Item(63, 11).Value = 3.3
$ws.Cells.Item(63, 12).Value = 3.25
$ws.Cells.Item(63, 13).Value = 2
$ws.Cells.Item(63, 14).Value = 3.5
$ws.Cells.Item(63, 15).Value = 3.25
$ws.Cells.Item(63, 16).Value = 1.909
$ws.Cells.Item(63, 17).Value = 0.5
$ws.Cells.Item(63, 18).Value = 1.775
$ws.Cells.Item(63, 19).Value = 2.025
$ws.Cells.Item(63, 20).Value = 2.25
$ws.Cells.Item(63, 23).Value = -1
$ws.Cells.Item(63, 24).Value = 2.25
$ws.Cells.Item(63, 26).Value = 0.7749999999999999
# Row 64
$ws.Cells.Item(64, 2).Value = 6221729
$ws.Cells.Item(64, 6).Value = "Kairat Almaty"
$ws.Cells.Item(64, 7).Value = "FC Astana"
$ws.Cells.Item(64, 8).Value = 1
$ws.Cells.Item(64, 10).Value = "H"
$ws.Cells.Item(64, 11).Value = 2.9
$ws.Cells.Item(64, 12).Value = 3.3
$ws.Cells.Item(64, 13).Value = 2.15
$ws.Cells.Item(64, 14).Value = 1.75
$ws.Cells.Item(64, 15).Value = 3.5
$ws.Cells.Item(64, 16).Value = 4
$ws.Cells.Item(64, 17).Value = -0.5
$ws.Cells.Item(64, 18).Value = 1.8
$ws.Cells.Item(64, 19).Value = 2
$ws.Cells.Item(64, 20).Value = 2.5
$ws.Cells.Item(64, 23).Value = 0.75
$ws.Cells.Item(64, 24).Value = -1
$ws.Cells.Item(64, 26).Value = 0.8
# Row 65
$ws.Cells.Item(65, 7).Value = "FK Aktobe"
# Row 70
$ws.Cells.Item(70, 6).Value = "Ordabasy"
$ws.Cells.Item(70, 7).Value = "FK Aktobe"
# Row 73
$ws.Cells.Item(73, 7).Value = "Ordabasy"
# Row 79
$ws.Cells.Item(79, 6).Value = "FK Aktobe"
# Row 82
$ws.Cells.Item(82, 6).Value = "FK Aktobe"
# Row 84
$ws.Cells.Item(84, 6).Value = "Ordabasy"
# Row 87
$ws.Cells.Item(87, 7).Value = "FK Aktobe"
# Row 95
$ws.Cells.Item(95, 7).Value = "Ordabasy"
# Row 97
$ws.Cells.Item(97, 6).Value = "FK Aktobe"
# Row 98
$ws.Cells.Item(98, 2).Value = 6221815
$ws.Cells.Item(98, 6).Value = "FK Atyrau"
$ws.Cells.Item(98, 7).Value = "Kairat Almaty"
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = "D"
$ws.Cells.Item(98, 11).Value = 3
$ws.Cells.Item(98, 12).Value = 3
$ws.Cells.Item(98, 13).Value = 2.25
$ws.Cells.Item(98, 14).Value = 3.1
$ws.Cells.Item(98, 15).Value = 3.1
$ws.Cells.Item(98, 16).Value = 2.15
$ws.Cells.Item(98, 17).Value = 0.25
$ws.Cells.Item(98, 18).Value = 1.85
$ws.Cells.Item(98, 19).Value = 1.95
$ws.Cells.Item(98, 20).Value = 2.25
$ws.Cells.Item(98, 21).Value = 1.8
$ws.Cells.Item(98, 22).Value = 2
$ws.Cells.Item(98, 24).Value = 2.1
$ws.Cells.Item(98, 25).Value = -1
$ws.Cells.Item(98, 26).Value = 0.425
$ws.Cells.Item(98, 27).Value = -0.5
$ws.Cells.Item(98, 28).Value = -1
$ws.Cells.Item(98, 29).Value = 1
# Row 99
$ws.Cells.Item(99, 2).Value = 6221753
$ws.Cells.Item(99, 6).Value = "FK Aksu"
$ws.Cells.Item(99, 7).Value = "Tobol Kostanay"
$ws.Cells.Item(99, 9).Value = 3
$ws.Cells.Item(99, 11).Value = 2.75
$ws.Cells.Item(99, 12).Value = 3.1
$ws.Cells.Item(99, 13).Value = 2.375
$ws.Cells.Item(99, 14).Value = 2.625
$ws.Cells.Item(99, 16).Value = 2.45
$ws.Cells.Item(99, 17).Value = 0
$ws.Cells.Item(99, 18).Value = 2
$ws.Cells.Item(99, 19).Value = 1.8
$ws.Cells.Item(99, 20).Value = 2.5
$ws.Cells.Item(99, 21).Value = 1.9
$ws.Cells.Item(99, 22).Value = 1.9
$ws.Cells.Item(99, 25).Value = 1.45
$ws.Cells.Item(99, 27).Value = 0.8
$ws.Cells.Item(99, 28).Value = 0.8999999999999999
$ws.Cells.Item(99, 29).Value = -1
# Row 100
$ws.Cells.Item(100, 2).Value = 6221752
$ws.Cells.Item(100, 6).Value = "FK Kyzylzhar"
$ws.Cells.Item(100, 7).Value = "Kaisar Kyzylorda"
$ws.Cells.Item(100, 9).Value = 1
$ws.Cells.Item(100, 10).Value = "A"
$ws.Cells.Item(100, 11).Value = 1.833
$ws.Cells.Item(100, 12).Value = 3.2
$ws.Cells.Item(100, 13).Value = 4
$ws.Cells.Item(100, 14).Value = 1.85
$ws.Cells.Item(100, 15).Value = 3.2
$ws.Cells.Item(100, 16).Value = 4
$ws.Cells.Item(100, 17).Value = -0.5
$ws.Cells.Item(100, 18).Value = 1.9
$ws.Cells.Item(100, 19).Value = 1.9
$ws.Cells.Item(100, 20).Value = 2
$ws.Cells.Item(100, 21).Value = 1.775
$ws.Cells.Item(100, 22).Value = 2.025
$ws.Cells.Item(100, 24).Value = -1
$ws.Cells.Item(100, 25).Value = 3
$ws.Cells.Item(100, 26).Value = -1
$ws.Cells.Item(100, 27).Value = 0.8999999999999999
$ws.Cells.Item(100, 29).Value = 1.025
# Row 101
$ws.Cells.Item(101, 2).Value = 6221814
$ws.Cells.Item(101, 6).Value = "Okzhetpes Kokshetau"
$ws.Cells.Item(101, 7).Value = "FK Maktaaral"
$ws.Cells.Item(101, 11).Value = 2.3
$ws.Cells.Item(101, 12).Value = 3.1
$ws.Cells.Item(101, 13).Value = 2.8
$ws.Cells.Item(101, 14).Value = 2.3
$ws.Cells.Item(101, 15).Value = 3.1
$ws.Cells.Item(101, 16).Value = 2.8
$ws.Cells.Item(101, 17).Value = 0
$ws.Cells.Item(101, 18).Value = 1.75
$ws.Cells.Item(101, 19).Value = 2.05
$ws.Cells.Item(101, 21).Value = 1.875
$ws.Cells.Item(101, 22).Value = 1.925
$ws.Cells.Item(101, 24).Value = 2.1
$ws.Cells.Item(101, 26).Value = 0
$ws.Cells.Item(101, 27).Value = -0
$ws.Cells.Item(101, 29).Value = 0.4625
# Row 102
$ws.Cells.Item(102, 2).Value = 6221754
$ws.Cells.Item(102, 6).Value = "Shakhter Karagandy"
$ws.Cells.Item(102, 7).Value = "FC Astana"
$ws.Cells.Item(102, 8).Value = 1
$ws.Cells.Item(102, 9).Value = 1
$ws.Cells.Item(102, 11).Value = 3.6
$ws.Cells.Item(102, 13).Value = 1.8
$ws.Cells.Item(102, 14).Value = 5
$ws.Cells.Item(102, 15).Value = 1.4
$ws.Cells.Item(102, 16).Value = 5
$ws.Cells.Item(102, 17).Value = 0.25
$ws.Cells.Item(102, 18).Value = 1.7
$ws.Cells.Item(102, 19).Value = 2.1
$ws.Cells.Item(102, 20).Value = 2.25
$ws.Cells.Item(102, 21).Value = 1.9
$ws.Cells.Item(102, 22).Value = 1.9
$ws.Cells.Item(102, 24).Value = 0.3999999999999999
$ws.Cells.Item(102, 26).Value = 0.35
$ws.Cells.Item(102, 27).Value = -0.5
$ws.Cells.Item(102, 28).Value = -0.5
$ws.Cells.Item(102, 29).Value = 0.45
# Row 103
$ws.Cells.Item(103, 2).Value = 6221755
$ws.Cells.Item(103, 6).Value = "Ordabasy"
$ws.Cells.Item(103, 7).Value = "Zhetysu"
$ws.Cells.Item(103, 8).Value = 2
$ws.Cells.Item(103, 9).Value = 2
$ws.Cells.Item(103, 11).Value = 1.727
$ws.Cells.Item(103, 12).Value = 3.5
$ws.Cells.Item(103, 13).Value = 4
$ws.Cells.Item(103, 14).Value = 1.444
$ws.Cells.Item(103, 15).Value = 4
$ws.Cells.Item(103, 16).Value = 6
$ws.Cells.Item(103, 17).Value = -1.25
$ws.Cells.Item(103, 18).Value = 1.975
$ws.Cells.Item(103, 19).Value = 1.825
$ws.Cells.Item(103, 20).Value = 2.75
$ws.Cells.Item(103, 21).Value = 1.8
$ws.Cells.Item(103, 22).Value = 2
$ws.Cells.Item(103, 24).Value = 3
$ws.Cells.Item(103, 26).Value = -1
$ws.Cells.Item(103, 27).Value = 0.825
$ws.Cells.Item(103, 28).Value = 0.8
$ws.Cells.Item(103, 29).Value = -1
# Row 104
$ws.Cells.Item(104, 7).Value = "FK Aktobe"
# Row 106
$ws.Cells.Item(106, 7).Value = "Ordabasy"
# Row 107
$ws.Cells.Item(107, 2).Value = 7874784
$ws.Cells.Item(107, 6).Value = "FC Elimai Semey"
$ws.Cells.Item(107, 7).Value = "FK Atyrau"
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 9).Value = 0
$ws.Cells.Item(107, 10).Value = "D"
$ws.Cells.Item(107, 11).Value = 2.45
$ws.Cells.Item(107, 12).Value = 3
$ws.Cells.Item(107, 13).Value = 2.7
$ws.Cells.Item(107, 14).Value = 2.45
$ws.Cells.Item(107, 15).Value = 3
$ws.Cells.Item(107, 16).Value = 2.7
$ws.Cells.Item(107, 17).Value = 0
$ws.Cells.Item(107, 18).Value = 1.8
$ws.Cells.Item(107, 19).Value = 2
$ws.Cells.Item(107, 20).Value = 2.25
$ws.Cells.Item(107, 21).Value = 2
$ws.Cells.Item(107, 22).Value = 1.8
$ws.Cells.Item(107, 23).Value = -1
$ws.Cells.Item(107, 24).Value = 2
$ws.Cells.Item(107, 26).Value = 0
$ws.Cells.Item(107, 27).Value = -0
$ws.Cells.Item(107, 28).Value = -1
$ws.Cells.Item(107, 29).Value = 0.8
# Row 108
$ws.Cells.Item(108, 2).Value = 7874783
$ws.Cells.Item(108, 6).Value = "Kairat Almaty"
$ws.Cells.Item(108, 7).Value = "FK Kyzylzhar"
$ws.Cells.Item(108, 8).Value = 2
$ws.Cells.Item(108, 9).Value = 1
$ws.Cells.Item(108, 10).Value = "H"
$ws.Cells.Item(108, 11).Value = 1.95
$ws.Cells.Item(108, 12).Value = 3.2
$ws.Cells.Item(108, 13).Value = 3.5
$ws.Cells.Item(108, 14).Value = 1.666
$ws.Cells.Item(108, 15).Value = 3.5
$ws.Cells.Item(108, 16).Value = 4.5
$ws.Cells.Item(108, 17).Value = -0.75
$ws.Cells.Item(108, 18).Value = 1.9
$ws.Cells.Item(108, 19).Value = 1.9
$ws.Cells.Item(108, 20).Value = 2.5
$ws.Cells.Item(108, 21).Value = 1.95
$ws.Cells.Item(108, 22).Value = 1.75
$ws.Cells.Item(108, 23).Value = 0.6659999999999999
$ws.Cells.Item(108, 24).Value = -1
$ws.Cells.Item(108, 26).Value = 0.45
$ws.Cells.Item(108, 27).Value = -0.5
$ws.Cells.Item(108, 28).Value = 0.95
$ws.Cells.Item(108, 29).Value = -1
# Row 109
$ws.Cells.Item(109, 6).Value = "FK Aktobe"
# Row 113
$ws.Cells.Item(113, 7).Value = "FK Aktobe"
# Row 119
$ws.Cells.Item(119, 2).Value = 7873759
$ws.Cells.Item(119, 6).Value = "Ordabasy"
$ws.Cells.Item(119, 7).Value = "FK Zhenys"
$ws.Cells.Item(119, 8).Value = 3
$ws.Cells.Item(119, 10).Value = "H"
$ws.Cells.Item(119, 11).Value = 1.25
$ws.Cells.Item(119, 12).Value = 5.75
$ws.Cells.Item(119, 13).Value = 7
$ws.Cells.Item(119, 14).Value = 1.444
$ws.Cells.Item(119, 15).Value = 4.75
$ws.Cells.Item(119, 16).Value = 4.75
$ws.Cells.Item(119, 17).Value = -1.25
$ws.Cells.Item(119, 18).Value = 1.95
$ws.Cells.Item(119, 19).Value = 1.85
$ws.Cells.Item(119, 20).Value = 2.75
$ws.Cells.Item(119, 21).Value = 1.9
$ws.Cells.Item(119, 22).Value = 1.9
$ws.Cells.Item(119, 23).Value = 0.444
$ws.Cells.Item(119, 24).Value = -1
$ws.Cells.Item(119, 26).Value = 0.95
$ws.Cells.Item(119, 27).Value = -1
$ws.Cells.Item(119, 28).Value = 0.45
$ws.Cells.Item(119, 29).Value = -0.5
# Row 120
$ws.Cells.Item(120, 2).Value = 7874795
$ws.Cells.Item(120, 6).Value = "FK Kyzylzhar"
$ws.Cells.Item(120, 7).Value = "Tobol Kostanay"
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 10).Value = "D"
$ws.Cells.Item(120, 11).Value = 2.2
$ws.Cells.Item(120, 12).Value = 3.1
$ws.Cells.Item(120, 13).Value = 3
$ws.Cells.Item(120, 14).Value = 2.625
$ws.Cells.Item(120, 15).Value = 3
$ws.Cells.Item(120, 16).Value = 2.55
$ws.Cells.Item(120, 17).Value = 0
$ws.Cells.Item(120, 18).Value = 1.9
$ws.Cells.Item(120, 19).Value = 1.9
$ws.Cells.Item(120, 20).Value = 2
$ws.Cells.Item(120, 21).Value = 1.95
$ws.Cells.Item(120, 22).Value = 1.85
$ws.Cells.Item(120, 23).Value = -1
$ws.Cells.Item(120, 24).Value = 2
$ws.Cells.Item(120, 26).Value = 0
$ws.Cells.Item(120, 27).Value = -0
$ws.Cells.Item(120, 28).Value = -1
$ws.Cells.Item(120, 29).Value = 0.8500000000000001
# Row 122
$ws.Cells.Item(122, 6).Value = "FK Aktobe"
# Row 126
$ws.Cells.Item(126, 6).Value = "Ordabasy"
# Row 128
$ws.Cells.Item(128, 7).Value = "FK Aktobe"
